$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.126.76"
$ws.Range("E2").Value = "  +0.99%  "

$ws.Range("D3").Value = "2.283.97"
$ws.Range("E3").Value = "  +3.59%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "273.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +5.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.66"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +11.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.627"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.45%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.632"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +6.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.14"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0940"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.66%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.11"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +9.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.57"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +8.14%  "

$ws.Range("D15").Value = "2.626.12"
$ws.Range("E15").Value = "  +3.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.836"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.97%  "

$ws.Range("D17").Value = "2.262.56"
$ws.Range("E17").Value = "  +1.08%  "

$ws.Range("D18").Value = "44.153.69"
$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("E19").Value = "  +2.63%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.07%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.18"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +14.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.30"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.64%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "237.19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.47"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +7.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +11.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.54"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.05%  "

$ws.Range("E29").Value = "  -5.31%  "

$ws.Range("E30").Value = "  -0.37%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.15"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +8.71%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0914"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.23%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.65"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +5.75%  "

$ws.Range("E35").Value = "  +1.69%  "

$ws.Range("E36").Value = "  +3.47%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0355"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.77%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.42"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.49"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +22.18%  "

$ws.Range("E40").Value = "  +26.09%  "

$ws.Range("E41").Value = "  +6.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.47"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.50"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "62.11"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.102"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +4.86%  "

$ws.Range("E46").Value = "  +3.21%  "

$ws.Range("E47").Value = "  +7.69%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("E49").Value = "  +1.64%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.434"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("D51").Value = "2.511.52"
$ws.Range("E51").Value = "  +3.74%  "
